# Update codigo e testes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapa2a")

# Row 12 (Logaritmico / Linear) - fill in full test parameters
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = "0-1"
$ws.Range("J12").Value = "0-1"
$ws.Range("L12").Value = 5
$ws.Range("M12").Value = 90
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = "0-1"
$ws.Range("S12").Value = "0-1"
$ws.Range("T12").Value = "Não"
$ws.Range("U12").Value = "Sim, no fim do mapa"
$ws.Range("V12").Value = 1
$ws.Range("W12").Value = 11

# Row 13 (Logaritmico / Gaussiano) - partial test parameters
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 30
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = "0-1"
$ws.Range("J13").Value = "0-1"
$ws.Range("L13").Value = 5
$ws.Range("M13").Value = 90
$ws.Range("N13").Value = 10

# Row 14 (Logaritmico / Logaritmico) - full test parameters
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "0-1"
$ws.Range("J14").Value = "0-1"
$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 90
$ws.Range("N14").Value = 10
$ws.Range("S14").Value = "0-2,3"

# New column (O) for the "Strenght Factor" field of the ResourceDetector block,
# plus a custom width for it.
$ws.Range("O5").Value = "Strenght Factor"
$ws.Columns.Item(15).ColumnWidth = 12.3

$ws.Range("O14").Value = 0.02
$ws.Range("T14").Value = "Sim"
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("U14").Value = "Não"
$ws.Range("V14").Value = 1
$ws.Range("W14").Value = 10
$ws.Range("R14").Value = "0-0,3"

# Leave the same selection state captured in the authored workbook.
$ws.Range("S15").Select()
